# Lecture title slide: bump the lecture number from 06 to 08
# (slide 1, "Text Placeholder 1" shape, second run of the first paragraph)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Only touch the " 06" run so the other run ("HWSW Flutter") stays untouched
# and the resulting XML keeps the same two-run paragraph layout.
$full = $tr.Text
$needle = " 06"
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $run = $tr.Characters($idx + 1, $needle.Length)
    $run.Text = " 08"
} else {
    # Fallback: if the expected text isn't found, just replace whole text.
    $tr.Text = $full -replace "06", "08"
}
